$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.197.72"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "1.843.27"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6878"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3015"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07661"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.93%  "

$ws.Range("D12").Value = "1.838.77"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.067"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6845"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.187"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.21%  "

$ws.Range("D17").Value = "29.193.03"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008183"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.89%  "

$ws.Range("D19").Value = "2.083.47"
$ws.Range("E19").Value = "  -0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.24%  "

$ws.Range("E21").Value = "  -1.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.413"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1453"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.785"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.517"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.285"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.150"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.197"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05248"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7612"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.855"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.136"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("D38").Value = "1.305.81"
$ws.Range("E38").Value = "  -0.88%  "

$ws.Range("E39").Value = "  -1.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.725"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9345"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.983"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "105.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9989"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").Value = "1.986.60"
$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5191"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000122"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.80%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.533"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.774"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "
